# daily auto push: 2026-02-05 14:10 UTC
# Inserts two new readings for 2026/02/05 (18:00 and 22:00) right after
# the existing 2026/02/05 rows, pushing every subsequent row down by 2
# (dimension grows from A1:D827 to A1:D829).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 786..827 down to 788..829, opening up two blank rows at 786:787.
$ws.Rows("786:787").Insert()

# Column A holds dates stored as plain text (e.g. "2026/02/05"), not real
# date serials, in this workbook. Force text formatting on column A before
# writing so Excel doesn't auto-coerce the string into a date value.
$ws.Range("A786:A787").NumberFormat = "@"

$ws.Range("A786").Value = "2026/02/05"
$ws.Range("B786").Value = "木"
$ws.Range("C786").Value = 18
$ws.Range("D786").Value = 201

$ws.Range("A787").Value = "2026/02/05"
$ws.Range("B787").Value = "木"
$ws.Range("C787").Value = 22
$ws.Range("D787").Value = 201

# Drop the explicit text format again so the new cells end up styleless,
# matching every other data row in the sheet (only the header row carries
# an explicit style).
$ws.Range("A786:D787").ClearFormats()
